$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7110151648521423
$ws.Range("B1").Value = 1.459963321685791
$ws.Range("C1").Value = 4.069529056549072
$ws.Range("D1").Value = 2.582533836364746
$ws.Range("E1").Value = 0.5626941919326782
